$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "Resolving-Mac"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.977367
$ws.Range("H2").Value2 = 2.932101
$ws.Range("I2").Value2 = 0.3701077125291157
$ws.Range("J2").Value2 = 0.3701077125291156
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.201254
$ws.Range("N2").Value2 = 0.603762
$ws.Range("Q2").Value2 = 0.196699018218
$ws.Range("R2").Value2 = 1.770291163962
$ws.Range("S2").Value2 = 0.3701077125291157
$ws.Range("T2").Value2 = 0.3701077125291156

# Row 3
$ws.Range("D3").Value2 = "Resolving-Mac"
$ws.Range("I3").Value2 = 0.1012017862004817
$ws.Range("J3").Value2 = 0.1012017862004817
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.201254
$ws.Range("N3").Value2 = 0.603762
$ws.Range("Q3").Value2 = 0.0537851315
$ws.Range("R3").Value2 = 0.4840661835
$ws.Range("S3").Value2 = 0.1012017862004817
$ws.Range("T3").Value2 = 0.1012017862004817

# Row 4
$ws.Range("D4").Value2 = "Resolving-Mac"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.4032770000000001
$ws.Range("H4").Value2 = 1.209831
$ws.Range("I4").Value2 = 0.1527122646719238
$ws.Range("J4").Value2 = 0.1527122646719238
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 0.201254
$ws.Range("N4").Value2 = 0.603762
$ws.Range("Q4").Value2 = 0.08116110935800001
$ws.Range("R4").Value2 = 0.7304499842220001
$ws.Range("S4").Value2 = 0.1527122646719238
$ws.Range("T4").Value2 = 0.1527122646719238

# Row 5
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.9928696666666667
$ws.Range("H5").Value2 = 2.978609
$ws.Range("I5").Value2 = 0.3759782365984789
$ws.Range("J5").Value2 = 0.375978236598479
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.201254
$ws.Range("N5").Value2 = 0.603762
$ws.Range("Q5").Value2 = 0.1998189918953333
$ws.Range("R5").Value2 = 1.798370927058
$ws.Range("S5").Value2 = 0.3759782365984789
$ws.Range("T5").Value2 = 0.375978236598479
